$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1. Insert 4 new rows at the top of the data (row 2), pushing existing data + formatting down
$ws.Range("A2:E5").EntireRow.Insert()

# 2. Copy formatting from row 6 (first of the shifted rows, which has the "normal" row style)
#    into the newly inserted blank rows 2-5
$ws.Range("A6:E6").Copy()
$ws.Range("A2:E5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Resize the table so it covers the new range
$tbl.Resize($ws.Range("A1:E20"))

# 4. Populate the Postcode/Date/News-Link cells for the new rows. Fill them in the same
#    order the new URLs appear in the target shared-string table (14-October,
#    16-october, 17-october) so the shared-string indices line up with the source data.
$ws.Cells.Item(5,1).Value = 44115
$ws.Cells.Item(5,2).Value = 3173
$ws.Cells.Item(5,3).Value = "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-14-October-2020"

$ws.Cells.Item(3,1).Value = 44117
$ws.Cells.Item(3,2).Value = 3047
$ws.Cells.Item(3,3).Value = "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-16-october-2020"

$ws.Cells.Item(4,1).Value = 44117
$ws.Cells.Item(4,2).Value = 3977
$ws.Cells.Item(4,3).Value = "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-16-october-2020"

$ws.Cells.Item(2,1).Value = 44118
$ws.Cells.Item(2,2).Value = 3128
$ws.Cells.Item(2,3).Value = "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020"

# 5. Re-apply the calculated-column formulas to every data row so the formula engine
#    re-resolves the structured references cleanly after the row insertion (the inserted
#    rows shifted formulas without re-binding them, leaving stale cached error values).
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r,4).Formula = "=VIC_Mystery_cases[[#This Row],[Date]]"
    $ws.Cells.Item($r,5).Formula = "=VIC_Mystery_cases[[#This Row],[Date]]+14"
}

# 6. Add hyperlink to C2 (pointing at the URL already present as its text)
$ws.Hyperlinks.Add($ws.Cells.Item(2,3), "https://www.dhhs.vic.gov.au/coronavirus-update-victoria-17-october-2020")

# 7. Re-apply the plain cell format to C2 since adding a hyperlink auto-applies the
#    "Hyperlink" style; the source data keeps the plain News-Link style for every row.
$ws.Range("C6").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 8. Set the active selection to C2, matching the saved view state
$ws.Range("C2").Select()

$wb.Save()
